$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 (kNN, SVM, LR, NB) with refreshed timing values ---

$row2 = @(0.01579561233520508, 0.02395977973937988, 0.009798669815063476, 0.01614923477172852, 0, 0.06241722106933593, 0.01784634590148926, 0.0213188648223877, 0.006253767013549805, 0.02848620414733887, 0.004152059555053711, 0.0256159782409668)
$row3 = @(0.0847860336303711, 0.0280029296875, 0.01257224082946777, 0.0115849494934082, 0.003669643402099609, 0.01166410446166992, 0.09928703308105469, 0.03016805648803711, 0.06685857772827149, 0.02468986511230469, 0.02614903450012207, 0.005408763885498047)
$row4 = @(0.02529869079589844, 0.01283197402954101, 0.01043562889099121, 0.01046676635742188, 0.08207383155822753, 0.008400058746337891, 0.02034802436828613, 0.01489152908325195, 0.01723241806030273, 0.0124995231628418, 0.03416523933410644, 0.006871843338012695)
$row5 = @(0.01709966659545898, 0.01573925018310547, 0.01494297981262207, 0.01405134201049805, 0, 0, 0.01505398750305176, 0.01183586120605469, 0.009855222702026368, 0.01453518867492676, 0, 0)

for ($i = 0; $i -lt $row2.Count; $i++) { $ws.Cells.Item(2, $i + 2).Value = $row2[$i] }
for ($i = 0; $i -lt $row3.Count; $i++) { $ws.Cells.Item(3, $i + 2).Value = $row3[$i] }
for ($i = 0; $i -lt $row4.Count; $i++) { $ws.Cells.Item(4, $i + 2).Value = $row4[$i] }
for ($i = 0; $i -lt $row5.Count; $i++) { $ws.Cells.Item(5, $i + 2).Value = $row5[$i] }

# --- Insert a new row 6 for "RF" (random forest), pushing the old "Ensemble" row down to row 7 ---

$ws.Rows.Item(6).Insert()

$ws.Cells.Item(6, 1).Value = "RF"
$ws.Cells.Item(6, 1).Font.Bold = $true
$ws.Cells.Item(6, 1).HorizontalAlignment = -4108
$ws.Cells.Item(6, 1).VerticalAlignment = -4160
$ws.Cells.Item(6, 1).Borders.LineStyle = 1

$rowRF = @(0.4757477283477783, 0.02020163536071777, 0.395454740524292, 0.01656913757324219, 0.8692080497741699, 0.0202143669128418, 0.239093542098999, 0.0164334774017334, 0.4434679985046387, 0.0183751106262207, 0.3708630561828613, 0.01582779884338379)
for ($i = 0; $i -lt $rowRF.Count; $i++) { $ws.Cells.Item(6, $i + 2).Value = $rowRF[$i] }

# --- Update the (now shifted) "Ensemble" row (row 7) with refreshed timing values ---

$ws.Cells.Item(7, 1).Value = "Ensemble"

$rowEnsemble = @(0.5574260711669922, 0.06136651039123535, 0.275451135635376, 0.03563728332519531, 0.6735668182373047, 0.03507504463195801, 0.9228752136230469, 0.09766650199890137, 0.3353415966033936, 0.04235577583312988, 0.5968525409698486, 0.03218626976013184)
for ($i = 0; $i -lt $rowEnsemble.Count; $i++) { $ws.Cells.Item(7, $i + 2).Value = $rowEnsemble[$i] }
